$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E
$ws.Range("E1").Value = "along"

# Updated data rows (row -> A label stays same, B/C/D updated, E new column)
$data = @(
    @{ Row = 2;  A = "<b>All</b>";       B = 0.0461285657788434;   C = 0.019797491413917;    D = 0.0724596401437698 },
    @{ Row = 3;  A = "<b>Europe</b>";    B = 0.0486706582010774;   C = 0.0110988437647862;   D = 0.0862424726373686 },
    @{ Row = 4;  A = "France";           B = 0.0249390799013659;   C = -0.0534596316251868;  D = 0.103337791427918 },
    @{ Row = 5;  A = "Germany";          B = 0.0772972876891449;   C = -0.0298528523416173;  D = 0.184447427719907 },
    @{ Row = 6;  A = "Italy";            B = 0.111299894074955;    C = -0.00805843128584366; D = 0.230658219435754 },
    @{ Row = 7;  A = "Poland";           B = 0.0336828202518876;   C = -0.0671457193603415;  D = 0.134511359864117 },
    @{ Row = 8;  A = "Spain";            B = 0.190567503358391;    C = 0.0828294959941132;   D = 0.298305510722669 },
    @{ Row = 9;  A = "United Kingdom";   B = -0.018749281740621;   C = -0.0994201362474925;  D = 0.0619215727662506 },
    @{ Row = 10; A = "Switzerland";      B = -0.00156935444217611; C = -0.168767094442469;   D = 0.165628385558117 },
    @{ Row = 11; A = "Japan";            B = 0.0255434199687464;   C = -0.0137594871745128;  D = 0.0648463271120057 },
    @{ Row = 12; A = "USA";              B = 0.104120455771318;    C = 0.0157757410245066;   D = 0.192465170518129 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.A
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = "millionaire_tax_in_programTRUE"
}
